# "Add files via upload" - refreshed Poker "Person Info" stats for the new
# month (sharePoint folder moved from "12 - Dec" to "02 - Feb"). Sessions
# played (NoOfSessions / NoOfYears columns), cumulative points and the
# average session points were updated for everyone who played during the
# period; Maisy also picked up a new best hand (Royal Flush) which replaces
# her previous best (Straight Flush) along with the date it was achieved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maisy (row 2)
$ws.Range("D2").Value = 125
$ws.Range("E2").Value = 14
$ws.Range("I2").Value = 451
$ws.Range("J2").Value = 3.61
$ws.Range("L2").Value = "01. Royal Flush"
$ws.Range("M2").Value = "Ace,10,Jack,Queen,King"
$ws.Range("N2").Value = 20.02

# Mark (row 3)
$ws.Range("D3").Value = 144
$ws.Range("E3").Value = 16
$ws.Range("I3").Value = 523
$ws.Range("J3").Value = 3.63

# Matt (row 4)
$ws.Range("D4").Value = 194
$ws.Range("E4").Value = 20
$ws.Range("I4").Value = 724
$ws.Range("J4").Value = 3.73

# Pepe (row 5)
$ws.Range("D5").Value = 106
$ws.Range("E5").Value = 13
$ws.Range("I5").Value = 383
$ws.Range("J5").Value = 3.61

# Prashant (row 6)
$ws.Range("D6").Value = 41
$ws.Range("E6").Value = 6
$ws.Range("I6").Value = 148
$ws.Range("J6").Value = 3.61

# Richard (row 7)
$ws.Range("D7").Value = 145
$ws.Range("E7").Value = 20
$ws.Range("I7").Value = 619
$ws.Range("J7").Value = 4.2699999999999996

# Jon (row 8)
$ws.Range("D8").Value = 199
$ws.Range("E8").Value = 20
$ws.Range("I8").Value = 702
$ws.Range("J8").Value = 3.53

# Alex (row 9)
$ws.Range("D9").Value = 86
$ws.Range("E9").Value = 15
$ws.Range("J9").Value = 3.42

# Andy (row 10)
$ws.Range("D10").Value = 202
$ws.Range("E10").Value = 20
$ws.Range("I10").Value = 833
$ws.Range("J10").Value = 4.12

# Anthony (row 11)
$ws.Range("D11").Value = 123
$ws.Range("E11").Value = 14
$ws.Range("I11").Value = 484
$ws.Range("J11").Value = 3.93
